# #3473 swapped out two properties
# Update two property square-footage/id values and adjust the
# selected range on the "BPS Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap out the two property values
$ws.Range("L3").Value = 227440.2
$ws.Range("B6").Value = 22482006
$ws.Range("B10").Value = 22482007

# Update the sheet's selection/scroll position (also clears the old
# topLeftCell scroll override, restoring the default view)
$ws.Range("L2:L10").Select()
